$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally has 19 columns (A..S). Two "best mean z-score over 5
# residue window" columns (I and M) are being dropped entirely, and the
# leftover indexed "regex_match_stpos_in_hit" column (S) is replaced by two
# new output columns (Vertebrata/Tetrapoda_regex_match_mean_zscore) appended
# at the end.
#
# Delete from right to left (S=19, M=13, I=9) so earlier deletions don't
# shift the column index of a later one.
$ws.Columns.Item(19).Delete()
$ws.Columns.Item(13).Delete()
$ws.Columns.Item(9).Delete()

# After the three deletions the sheet is A..P (16 cols). Append the two new
# header columns at Q/R, copying the bold/centered/bordered header format
# from the neighboring header cell (P1) so they share its style record.
$ws.Range("P1").Copy()
$ws.Range("Q1:R1").PasteSpecial(-4122)

$ws.Range("Q1").Value = "Vertebrata_regex_match_mean_zscore"
$ws.Range("R1").Value = "Tetrapoda_regex_match_mean_zscore"

# Only the row-4 record had a regex match, so only it gets the new scores.
$ws.Range("Q4").Value = -0.7300460456318796
$ws.Range("R4").Value = -0.3685849517941594

# The other data rows still carry a (blank) cell all the way out to column R
# -- force that by touching formatting on the trailing blank cells so they
# aren't trimmed from the row on save.
$ws.Range("Q2:R3").Style = $ws.Range("P2").Style
$ws.Range("Q5:R8").Style = $ws.Range("P2").Style
$ws.Range("P5").Style = $ws.Range("P2").Style
